$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Appointment ID values in column A to the new AP0x scheme
$ws.Range("A2").Value = "AP05"
$ws.Range("A3").Value = "AP02"
$ws.Range("A4").Value = "AP03"
$ws.Range("A5").Value = "AP01"

# The custom date format used by B5 switches from upper-case MM to lower-case mm
$ws.Range("B5").NumberFormat = "dd/mm/yy"

# Update the active selection shown when the sheet is next opened
$ws.Range("A9").Select() | Out-Null
